$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after the existing "AddCustomerTest" sheet
$afterSheet = $wb.Worksheets.Item("AddCustomerTest")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "OpenAccountTest"

# Header row
$newSheet.Range("A1").Value = "customer"
$newSheet.Range("B1").Value = "currency"

# Data row
$newSheet.Range("A2").Value = "Huyen Ha"
$newSheet.Range("B2").Value = "VND"

# Leave the cursor on B2 (last entered cell) and make the new sheet active/selected
[void]$newSheet.Range("B2").Select()
$newSheet.Activate()
